# Bitacora de iteraciones - agregar fila de la iteracion 3
# (actualizacion de glosario y alcance)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bitacora")

# --- Nueva fila 11: datos de la iteracion 3 ---

# A11 "Iteracion" se escribe como texto "3.0" (igual que 1.0 / 2.0 arriba)
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "3.0"
$ws.Range("A11").Style = "Normal"

$ws.Range("B11").Value = "H1 – Bicicletas"
$ws.Range("C11").Value = "feature/creacion_historias_usuario_iter_3"
$ws.Range("D11").Value = "creacion_de_historias"
$ws.Range("E11").Value = "andr4f"
$ws.Range("F11").Value = "andr4f, Angel Trillo, Yineth Avila"

$ws.Hyperlinks.Add($ws.Range("G11"), "https://github.com/andr4f/bici-go-bd/pull/16")
$ws.Range("G11").Value = "https://github.com/andr4f/bici-go-bd/pull/16"
$ws.Range("G10").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("H11").Value = "v0.1-iter3-creacion"
$ws.Range("I11").Value = "Creacion de 4 historias mas para completar 10, creacion de glosario y alcance de datos"

$ws.Range("J11").Value = "10/14/2025"
$ws.Range("K11").Value = "10/15/2025"
$ws.Range("J11").NumberFormat = $ws.Range("J10").NumberFormat
$ws.Range("K11").NumberFormat = $ws.Range("K10").NumberFormat

$ws.Range("L11").Value = "Completado"
$ws.Range("M11").Value = "documentacion/backlog_historias_usuario/"

# El nuevo comentario es el texto mas largo de la columna I -> ajustar ancho
$ws.Columns("I").ColumnWidth = 72.3

# Seleccion final tal como quedo guardada
$ws.Activate()
$ws.Range("D27").Select()
